# Auto-generated edit script for Nrg2-Erbb3.xlsx TPM data update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nrg2"
$ws.Range("C2").Value = "Erbb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1182943333333333
$ws.Range("H2").Value = 0.354883
$ws.Range("I2").Value = 0.07902913105657369
$ws.Range("J2").Value = 0.07902913105657366
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2023976666666667
$ws.Range("N2").Value = 0.6071930000000001
$ws.Range("O2").Value = 0.03663970451354832
$ws.Range("P2").Value = 0.03663970451354832
$ws.Range("Q2").Value = 0.02394249704655556
$ws.Range("R2").Value = 0.215482473419
$ws.Range("S2").Value = 0.002895604009875344
$ws.Range("T2").Value = 0.002895604009875344

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nrg2"
$ws.Range("C3").Value = "Erbb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1182943333333333
$ws.Range("H3").Value = 0.354883
$ws.Range("I3").Value = 0.07902913105657369
$ws.Range("J3").Value = 0.07902913105657366
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.4773683333333333
$ws.Range("N3").Value = 1.432105
$ws.Range("O3").Value = 0.08641717548188978
$ws.Range("P3").Value = 0.08641717548188979
$ws.Range("Q3").Value = 0.05646996874611111
$ws.Range("R3").Value = 0.508229718715
$ws.Range("S3").Value = 0.006829474286697194
$ws.Range("T3").Value = 0.006829474286697192

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Nrg2"
$ws.Range("C4").Value = "Erbb3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1182943333333333
$ws.Range("H4").Value = 0.354883
$ws.Range("I4").Value = 0.07902913105657369
$ws.Range("J4").Value = 0.07902913105657366
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.844232333333333
$ws.Range("N4").Value = 14.532697
$ws.Range("O4").Value = 0.876943120004562
$ws.Range("P4").Value = 0.876943120004562
$ws.Range("Q4").Value = 0.5730452343834445
$ws.Range("R4").Value = 5.157407109450999
$ws.Range("S4").Value = 0.06930405276000116
$ws.Range("T4").Value = 0.06930405276000114

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Nrg2"
$ws.Range("C5").Value = "Erbb3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.178245666666667
$ws.Range("H5").Value = 3.534737
$ws.Range("I5").Value = 0.7871529310322559
$ws.Range("J5").Value = 0.7871529310322558
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2023976666666667
$ws.Range("N5").Value = 0.6071930000000001
$ws.Range("O5").Value = 0.03663970451354832
$ws.Range("P5").Value = 0.03663970451354832
$ws.Range("Q5").Value = 0.2384741736934445
$ws.Range("R5").Value = 2.146267563241
$ws.Range("S5").Value = 0.02884105079999533
$ws.Range("T5").Value = 0.02884105079999533

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Nrg2"
$ws.Range("C6").Value = "Erbb3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.178245666666667
$ws.Range("H6").Value = 3.534737
$ws.Range("I6").Value = 0.7871529310322559
$ws.Range("J6").Value = 0.7871529310322558
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.4773683333333333
$ws.Range("N6").Value = 1.432105
$ws.Range("O6").Value = 0.08641717548188978
$ws.Range("P6").Value = 0.08641717548188979
$ws.Range("Q6").Value = 0.5624571701538889
$ws.Range("R6").Value = 5.062114531384999
$ws.Range("S6").Value = 0.06802353297209834
$ws.Range("T6").Value = 0.06802353297209834

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Nrg2"
$ws.Range("C7").Value = "Erbb3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.178245666666667
$ws.Range("H7").Value = 3.534737
$ws.Range("I7").Value = 0.7871529310322559
$ws.Range("J7").Value = 0.7871529310322558
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.844232333333333
$ws.Range("N7").Value = 14.532697
$ws.Range("O7").Value = 0.876943120004562
$ws.Range("P7").Value = 0.876943120004562
$ws.Range("Q7").Value = 5.707695755076554
$ws.Range("R7").Value = 51.36926179568899
$ws.Range("S7").Value = 0.6902883472601623
$ws.Range("T7").Value = 0.6902883472601622

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Nrg2"
$ws.Range("C8").Value = "Erbb3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.2003046666666667
$ws.Range("H8").Value = 0.6009139999999999
$ws.Range("I8").Value = 0.1338179379111705
$ws.Range("J8").Value = 0.1338179379111705
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2023976666666667
$ws.Range("N8").Value = 0.6071930000000001
$ws.Range("O8").Value = 0.03663970451354832
$ws.Range("P8").Value = 0.03663970451354832
$ws.Range("Q8").Value = 0.04054119715577778
$ws.Range("R8").Value = 0.364870774402
$ws.Range("S8").Value = 0.004903049703677642
$ws.Range("T8").Value = 0.004903049703677641

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Nrg2"
$ws.Range("C9").Value = "Erbb3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.2003046666666667
$ws.Range("H9").Value = 0.6009139999999999
$ws.Range("I9").Value = 0.1338179379111705
$ws.Range("J9").Value = 0.1338179379111705
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4773683333333333
$ws.Range("N9").Value = 1.432105
$ws.Range("O9").Value = 0.08641717548188978
$ws.Range("P9").Value = 0.08641717548188979
$ws.Range("Q9").Value = 0.09561910488555556
$ws.Range("R9").Value = 0.8605719439699999
$ws.Range("S9").Value = 0.01156416822309425
$ws.Range("T9").Value = 0.01156416822309425

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Nrg2"
$ws.Range("C10").Value = "Erbb3"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.2003046666666667
$ws.Range("H10").Value = 0.6009139999999999
$ws.Range("I10").Value = 0.1338179379111705
$ws.Range("J10").Value = 0.1338179379111705
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.844232333333333
$ws.Range("N10").Value = 14.532697
$ws.Range("O10").Value = 0.876943120004562
$ws.Range("P10").Value = 0.876943120004562
$ws.Range("Q10").Value = 0.9703223427842221
$ws.Range("R10").Value = 8.732901085057998
$ws.Range("S10").Value = 0.1173507199843986
$ws.Range("T10").Value = 0.1173507199843986
